# Append the new data row (2026/01/25) to the "ModCounts" sheet, mirroring
# the formatting of the preceding row (row 75).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ModCounts")

$newRow = 76
$prevRow = 75

# Write the new values. Column A holds a date formatted as plain text
# (e.g. "2026/01/25"), just like all the other rows in the sheet, so a
# leading apostrophe forces it to stay text instead of being auto-parsed
# into a date serial number.
$ws.Range("A" + $newRow).Value = "'2026/01/25"
$ws.Range("B" + $newRow).Value = "逃离鸭科夫"
$ws.Range("C" + $newRow).Value = 1159

# Match the style (center alignment, etc.) used by the rest of the table
# by copying the formatting from the previous row.
$ws.Range("A" + $prevRow + ":C" + $prevRow).Copy()
$ws.Range("A" + $newRow + ":C" + $newRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0
